$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 38 (pushes existing rows 38-51 down to 40-53)
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

# --- New row 38: weekly update, "Extra" bundled as $/bandeja 10 kilos ---
$ws.Cells.Item(38, 1).Value = 10
$ws.Cells.Item(38, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value = "La Araucanía"
$ws.Cells.Item(38, 4).Value = 44508
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 300000000
$ws.Cells.Item(38, 7).Value = "Espárragos"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 90
$ws.Cells.Item(38, 11).Value = 14000
$ws.Cells.Item(38, 12).Value = 14000
$ws.Cells.Item(38, 13).Value = 14000
$ws.Cells.Item(38, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(38, 15).Value = "Región del Maule"
$ws.Cells.Item(38, 16).Value = 1400
$ws.Cells.Item(38, 17).Value = 10
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# --- New row 39: weekly update, regular $/kilo pricing ---
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 44508
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 300000000
$ws.Cells.Item(39, 7).Value = "Espárragos"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 100
$ws.Cells.Item(39, 11).Value = 1400
$ws.Cells.Item(39, 12).Value = 1400
$ws.Cells.Item(39, 13).Value = 1400
$ws.Cells.Item(39, 14).Value = "$/kilo"
$ws.Cells.Item(39, 15).Value = "Región del Maule"
$ws.Cells.Item(39, 16).Value = 1400
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
